$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for the team win/loss/tie record columns (AC, AD, AE)
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the header formatting used by the other header cells (bold, thin
# border, centered horizontally, top-aligned vertically)
$headerFmt = $ws.Range("AB1")
$newHeaders = $ws.Range("AC1:AE1")
$newHeaders.Font.Bold = $headerFmt.Font.Bold
$newHeaders.HorizontalAlignment = $headerFmt.HorizontalAlignment
$newHeaders.VerticalAlignment = $headerFmt.VerticalAlignment
$newHeaders.Borders.LineStyle = $headerFmt.Borders.LineStyle

# Fill in the team's record for every player row (same W/L/T for the team)
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 29).Value = 78
    $ws.Cells.Item($r, 30).Value = 84
    $ws.Cells.Item($r, 31).Value = 0
}
